$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.689.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "'3.133.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'567.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "'148.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.33%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'3.131.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("D11").Value = "'6.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").Value = "'0.500"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.34%  "
$ws.Range("D13").Value = "'0.0000269"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +17.06%  "
$ws.Range("D14").Value = "'37.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.36%  "
$ws.Range("D15").Value = "'3.646.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "'64.806.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'7.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.52%  "
$ws.Range("D18").Value = "'3.140.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").Value = "'506.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.83%  "
$ws.Range("D21").Value = "'14.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.33%  "
$ws.Range("D22").Value = "'0.730"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.90%  "
$ws.Range("E23").Value = "  +10.06%  "
$ws.Range("D24").Value = "'7.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.22%  "
$ws.Range("D25").Value = "'84.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.05%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'2.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.41%  "
$ws.Range("D28").Value = "'8.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.22%  "
$ws.Range("D29").Value = "'2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.79%  "
$ws.Range("D30").Value = "'27.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.54%  "
$ws.Range("B31").Value = "Mantle"
$ws.Range("C31").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D31").Value = "'1.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.71%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.45%  "
$ws.Range("E34").Value = "  +11.60%  "
$ws.Range("D35").Value = "'6.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.80%  "
$ws.Range("D36").Value = "'55.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").Value = "'473.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.27%  "
$ws.Range("D38").Value = "'0.0860"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.32%  "
$ws.Range("D39").Value = "'0.0416"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("D40").Value = "'2.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("D41").Value = "'3.107.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.64%  "
$ws.Range("D42").Value = "'8.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.17%  "
$ws.Range("E43").Value = "  +4.86%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.69%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.288"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.65%  "
$ws.Range("D46").Value = "'29.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.11%  "
$ws.Range("D47").Value = "0.0₃0573"
$ws.Range("E47").Value = "  +11.91%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'0.115"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("D50").Value = "'2.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.16%  "
$ws.Range("D51").Value = "'122.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.59%  "
